# Apply the updates described in the diff to the "Inscricoes" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 7: Inscritos (E) 4 -> 5
$ws.Range("E7").Value = 5

# Row 17: Inscritos (E) 59 -> 60
$ws.Range("E17").Value = 60

# Row 20: Inscritos (E) 2 -> 3
$ws.Range("E20").Value = 3

# Row 35: Inscritos (E) 2 -> 3
$ws.Range("E35").Value = 3

# Row 36: Pagos (F) 17 -> 18, Inscricoes homologadas (H) 17 -> 18
$ws.Range("F36").Value = 18
$ws.Range("H36").Value = 18

# Row 38: Inscritos (E) 36 -> 37, Pagos (F) 8 -> 9, Inscricoes homologadas (H) 8 -> 9
$ws.Range("E38").Value = 37
$ws.Range("F38").Value = 9
$ws.Range("H38").Value = 9

# Row 39: Inscritos (E) 12 -> 13, Pagos (F) 6 -> 7, Inscricoes homologadas (H) 6 -> 7
$ws.Range("E39").Value = 13
$ws.Range("F39").Value = 7
$ws.Range("H39").Value = 7

# Row 79: Inscritos (E) 14 -> 15
$ws.Range("E79").Value = 15
